# edit.ps1 - applies the recorded change:
#   1. Swap the design theme applied to the deck's slide master from the
#      "Integral" (Red Violet) palette to the standard "Office Theme" palette
#      (this is what a user does via Design > Themes > Office Theme).
#   2. Re-apply a table style ("No Style, Table Grid" -
#      {34239943-07C3-43FD-A581-8ACC77DF4ED7}) to the three tables living on
#      slides 14, 15 and 16 (they previously used the deck's custom table
#      style {83CD74C5-4D8A-4C6C-BBBC-4DE51E11D50F}).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Theme colour swap: Integral/"Red Violet" -> "Office" colours.
# ---------------------------------------------------------------------
# PowerPoint's RGB COM property packs colour bytes as 0x00BBGGRR (like the
# Win32 COLORREF), so convert each RRGGBB hex triplet accordingly.
function New-BgrFromRgbHex([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Order matches ThemeColorScheme.Colors(1..12):
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6,
# hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = New-BgrFromRgbHex $officeThemeColors[$i - 1]
}

# ---------------------------------------------------------------------
# 2) Table style re-application on slides 14, 15 and 16.
# ---------------------------------------------------------------------
$newTableStyleId = "{34239943-07C3-43FD-A581-8ACC77DF4ED7}"
$tableSlideIndexes = @(14, 15, 16)

foreach ($slideIndex in $tableSlideIndexes) {
    $s = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId, $false)
        }
    }
}
